$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.454.58"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "3.015.69"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'584.17"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "'147.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.39%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").Value = "3.014.55"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -4.42%  "
$ws.Range("D11").Value = "'5.71"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "'0.442"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "'34.74"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.24%  "
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "3.513.31"
$ws.Range("E16").Value = "  -1.70%  "
$ws.Range("D17").Value = "'7.06"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").Value = "62.412.25"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").Value = "3.011.53"
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").Value = "'461.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").Value = "'13.91"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").Value = "'7.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "'2.29"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.78%  "
$ws.Range("D25").Value = "'79.98"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").Value = "'12.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.39%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'9.93"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").Value = "'7.16"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'2.10"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'27.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").Value = "0.0₃0786"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").Value = "'50.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'9.04"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("E41").Value = "  -10.95%  "
$ws.Range("D42").Value = "'416.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "'0.275"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.772.30"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0352"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "'37.96"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.85%  "
$ws.Range("D48").Value = "'129.10"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "'23.78"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.90%  "
